$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.948.41'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.265.53'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.15'
$ws.Range("E5").Value = '  +1.84%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.654'
$ws.Range("E6").Value = '  +4.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.50'
$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.451'
$ws.Range("E9").Value = '  +6.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0977'
$ws.Range("E10").Value = '  -2.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.23'
$ws.Range("E11").Value = '  +1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.75'
$ws.Range("E12").Value = '  +4.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  +1.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.601.07'
$ws.Range("E14").Value = '  -0.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.64'
$ws.Range("E15").Value = '  +0.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.17'
$ws.Range("E16").Value = '  +5.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.843'
$ws.Range("E17").Value = '  +4.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.262.25'
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.883.87'
$ws.Range("E19").Value = '  +0.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0986'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.08'
$ws.Range("E21").Value = '  +1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  +2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.55'
$ws.Range("E23").Value = '  +1.25%  '

$ws.Range("E24").Value = '  -0.25%  '

$ws.Range("E25").Value = '  -0.56%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.35'
$ws.Range("E27").Value = '  +20.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.94'
$ws.Range("E28").Value = '  +1.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.30'
$ws.Range("E29").Value = '  +9.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.05'
$ws.Range("E30").Value = '  +1.48%  '

$ws.Range("E31").Value = '  +0.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.44'
$ws.Range("E32").Value = '  +0.53%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.127'
$ws.Range("E33").Value = '  +4.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.98'
$ws.Range("E34").Value = '  +6.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0688'
$ws.Range("E35").Value = '  +0.39%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.97'
$ws.Range("E36").Value = '  -1.76%  '

$ws.Range("E37").Value = '  -2.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.45'
$ws.Range("E38").Value = '  -3.03%  '

$ws.Range("E39").Value = '  -0.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0256'

$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.80'
$ws.Range("E42").Value = '  +5.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000223'
$ws.Range("E43").Value = '  +5.98%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.83'
$ws.Range("E44").Value = '  +2.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.24'
$ws.Range("E45").Value = '  +1.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0951'
$ws.Range("E46").Value = '  -0.92%  '

$ws.Range("E47").Value = '  -0.28%  '

$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.456.44'
$ws.Range("E49").Value = '  -0.85%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.34'
$ws.Range("E50").Value = '  +0.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.96'
$ws.Range("E51").Value = '  -4.55%  '
